# Restore D8 ("Integer max" for rule R10 on the Rules sheet) from 11 to 13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D8").Value = 13
